# Update Alcachofa price records: re-sync dates, volumes, prices,
# variety/quality and origin fields per the weekly fruit/vegetable refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 44488
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 12000
$ws.Range("M2").Value = 11500
$ws.Range("P2").Value = 288
$ws.Range("D3").Value = 44425
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14500
$ws.Range("O3").Value = "Región del Maule"
$ws.Range("P3").Value = 362
$ws.Range("D4").Value = 44503
$ws.Range("J4").Value = 160
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 12000
$ws.Range("M4").Value = 11500
$ws.Range("P4").Value = 288
$ws.Range("D6").Value = 44484
$ws.Range("J6").Value = 120
$ws.Range("D7").Value = 44468
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 12000
$ws.Range("L7").Value = 13000
$ws.Range("M7").Value = 12500
$ws.Range("P7").Value = 312
$ws.Range("D8").Value = 44512
$ws.Range("D9").Value = 44432
$ws.Range("J9").Value = 120
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 14500
$ws.Range("P9").Value = 362
$ws.Range("D10").Value = 44496
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 12000
$ws.Range("M10").Value = 11500
$ws.Range("O10").Value = "Provincia del Elquí"
$ws.Range("P10").Value = 288
$ws.Range("D12").Value = 44453
$ws.Range("K12").Value = 12500
$ws.Range("L12").Value = 13000
$ws.Range("M12").Value = 12750
$ws.Range("O12").Value = "Provincia del Elquí"
$ws.Range("P12").Value = 319
$ws.Range("D13").Value = 44417
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 16000
$ws.Range("M13").Value = 15500
$ws.Range("P13").Value = 388
$ws.Range("D14").Value = 44495
$ws.Range("J14").Value = 120
$ws.Range("D15").Value = 44505
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 11000
$ws.Range("L15").Value = 12000
$ws.Range("M15").Value = 11500
$ws.Range("P15").Value = 288
$ws.Range("D16").Value = 44508
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 11000
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = 11500
$ws.Range("P16").Value = 288
$ws.Range("D18").Value = 44467
$ws.Range("O18").Value = "Provincia de Limarí"
$ws.Range("D19").Value = 44510
$ws.Range("D20").Value = 44475
$ws.Range("K20").Value = 11000
$ws.Range("L20").Value = 12000
$ws.Range("M20").Value = 11500
$ws.Range("P20").Value = 288
$ws.Range("D21").Value = 44482
$ws.Range("J21").Value = 120
$ws.Range("D22").Value = 44515
$ws.Range("J22").Value = 120
$ws.Range("K22").Value = 11000
$ws.Range("L22").Value = 12000
$ws.Range("M22").Value = 11500
$ws.Range("P22").Value = 288
$ws.Range("D23").Value = 44494
$ws.Range("D24").Value = 44454
$ws.Range("K24").Value = 13000
$ws.Range("L24").Value = 14000
$ws.Range("M24").Value = 13500
$ws.Range("P24").Value = 338
$ws.Range("D25").Value = 44420
$ws.Range("K25").Value = 13000
$ws.Range("L25").Value = 14000
$ws.Range("M25").Value = 13500
$ws.Range("P25").Value = 338
$ws.Range("D26").Value = 44516
$ws.Range("K26").Value = 11000
$ws.Range("L26").Value = 12000
$ws.Range("M26").Value = 11500
$ws.Range("O26").Value = "Provincia del Elquí"
$ws.Range("P26").Value = 288
$ws.Range("D27").Value = 44426
$ws.Range("K27").Value = 13000
$ws.Range("L27").Value = 14000
$ws.Range("M27").Value = 13500
$ws.Range("O27").Value = "Región del Maule"
$ws.Range("P27").Value = 338
$ws.Range("D28").Value = 44473
$ws.Range("J28").Value = 160
$ws.Range("K28").Value = 11000
$ws.Range("L28").Value = 12000
$ws.Range("M28").Value = 11500
$ws.Range("P28").Value = 288
$ws.Range("D29").Value = 44455
$ws.Range("J29").Value = 100
$ws.Range("K29").Value = 13000
$ws.Range("L29").Value = 14000
$ws.Range("M29").Value = 13500
$ws.Range("P29").Value = 338
$ws.Range("D30").Value = 44490
$ws.Range("J30").Value = 100
$ws.Range("D31").Value = 44446
$ws.Range("H31").Value = "Madrigal"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 160
$ws.Range("K31").Value = 12500
$ws.Range("L31").Value = 13000
$ws.Range("M31").Value = 12750
$ws.Range("P31").Value = 319
$ws.Range("D32").Value = 44489
$ws.Range("D33").Value = 44487
$ws.Range("J33").Value = 100
$ws.Range("D34").Value = 44399
$ws.Range("H34").Value = "Española"
$ws.Range("I34").Value = "Segunda"
$ws.Range("J34").Value = 120
$ws.Range("K34").Value = 15500
$ws.Range("L34").Value = 16000
$ws.Range("M34").Value = 15750
$ws.Range("P34").Value = 394
$ws.Range("D35").Value = 44498
$ws.Range("K35").Value = 10500
$ws.Range("L35").Value = 11000
$ws.Range("M35").Value = 10750
$ws.Range("P35").Value = 269
